$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  28"
$ws.Range("C9").Value = "Report Covering the Week  7/10/2023  Through  7/16/2023"

# --- Crime Complaints table updates ---
# Row 14
$ws.Range("N14").Value = -81.818181818181

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("C15").NumberFormat = "#,##0"
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = "#,##0"
$ws.Range("E15").Value = 0
$ws.Range("E15").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 6
$ws.Range("J15").Value = 7
$ws.Range("K15").Value = -14.285714285714
$ws.Range("L15").Value = -14.285714285714
$ws.Range("M15").Value = -25
$ws.Range("N15").Value = -62.5

# Row 16
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 4
$ws.Range("H16").Value = 175
$ws.Range("I16").Value = 57
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = 3.636363636363
$ws.Range("L16").Value = 62.857142857142
$ws.Range("M16").Value = -55.813953488372
$ws.Range("N16").Value = -84.759358288770

# Row 17
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 50
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = 72.727272727272
$ws.Range("I17").Value = 94
$ws.Range("J17").Value = 88
$ws.Range("K17").Value = 6.818181818181
$ws.Range("L17").Value = 6.818181818181
$ws.Range("M17").Value = 46.875
$ws.Range("N17").Value = -45.664739884393

# Row 18
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -75
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = -28.571428571428
$ws.Range("J18").Value = 54
$ws.Range("K18").Value = -5.555555555555
$ws.Range("M18").Value = -62.5
$ws.Range("N18").Value = -92.330827067669

# Row 19
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 61.111111111111
$ws.Range("I19").Value = 332
$ws.Range("J19").Value = 260
$ws.Range("K19").Value = 27.692307692307
$ws.Range("L19").Value = 80.434782608695
$ws.Range("M19").Value = 33.333333333333
$ws.Range("N19").Value = -7.262569832402

# Row 20
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("F20").Value = 18
$ws.Range("G20").Value = 12
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 64
$ws.Range("K20").Value = 4.6875
$ws.Range("L20").Value = 52.272727272727
$ws.Range("M20").Value = -21.176470588235
$ws.Range("N20").Value = -95.600787918581

# Row 21
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 31.818181818181
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 56.944444444444
$ws.Range("I21").Value = 609
$ws.Range("J21").Value = 528
$ws.Range("K21").Value = 15.340909090909
$ws.Range("L21").Value = 46.746987951807
$ws.Range("M21").Value = -9.239940387481
$ws.Range("N21").Value = -80.480769230769

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 500
$ws.Range("I23").Value = 13
$ws.Range("K23").Value = -35
$ws.Range("L23").Value = 30
$ws.Range("M23").Value = -23.529411764705

# Row 24
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = 3.333333333333
$ws.Range("F24").Value = 101
$ws.Range("G24").Value = 102
$ws.Range("H24").Value = -0.980392156862
$ws.Range("I24").Value = 656
$ws.Range("J24").Value = 590
$ws.Range("K24").Value = 11.186440677966
$ws.Range("L24").Value = 56.937799043062
$ws.Range("M24").Value = 34.426229508196

# Row 25
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = -71.428571428571
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = -4.761904761904
$ws.Range("I25").Value = 140
$ws.Range("J25").Value = 123
$ws.Range("K25").Value = 13.821138211382
$ws.Range("L25").Value = 12
$ws.Range("M25").Value = -24.731182795698

# Row 26
$ws.Range("C26").Value = 1
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("D26").Value = 1
$ws.Range("D26").NumberFormat = "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("E26").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("I26").Value = 14
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = 27.272727272727
$ws.Range("L26").Value = 16.666666666666

# Row 27
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 17
$ws.Range("J27").Value = 19
$ws.Range("K27").Value = -10.526315789473
$ws.Range("L27").Value = 41.666666666666

# Row 28
$ws.Range("L28").Value = -50
$ws.Range("N28").Value = -82.142857142857

# Row 29
$ws.Range("L29").Value = -50
$ws.Range("N29").Value = -84

